$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text format on every touched cell first, so that numeric-looking
# strings (e.g. "567.74", "0.158", "0.510") are stored verbatim as text
# instead of being auto-parsed into numbers -- matching the original
# inline-string (text) cell type used throughout this sheet.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("B33").NumberFormat = "@"
$ws.Range("C33").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("B34").NumberFormat = "@"
$ws.Range("C34").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("B42").NumberFormat = "@"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("B43").NumberFormat = "@"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("E51").NumberFormat = "@"

$ws.Range('D2').Value = '69.446.37'
$ws.Range('E2').Value = '  -0.25%  '
$ws.Range('D3').Value = '2.491.07'
$ws.Range('E3').Value = '  -0.67%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = '567.74'
$ws.Range('E5').Value = '  -1.00%  '
$ws.Range('D6').Value = '164.16'
$ws.Range('E6').Value = '  -0.45%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('E8').Value = '  -1.27%  '
$ws.Range('D9').Value = '2.490.30'
$ws.Range('E9').Value = '  -0.68%  '
$ws.Range('D10').Value = '0.158'
$ws.Range('E10').Value = '  -1.73%  '
$ws.Range('E11').Value = '  -0.09%  '
$ws.Range('D12').Value = '0.354'
$ws.Range('E12').Value = '  -0.20%  '
$ws.Range('D13').Value = '4.89'
$ws.Range('E13').Value = '  -0.75%  '
$ws.Range('D14').Value = '2.945.85'
$ws.Range('E14').Value = '  -0.63%  '
$ws.Range('D15').Value = '69.333.99'
$ws.Range('E15').Value = '  -0.28%  '
$ws.Range('E16').Value = '  -1.14%  '
$ws.Range('D17').Value = '24.28'
$ws.Range('E17').Value = '  -2.69%  '
$ws.Range('D18').Value = '2.506.33'
$ws.Range('E18').Value = '  -0.24%  '
$ws.Range('D19').Value = '11.19'
$ws.Range('E19').Value = '  -1.63%  '
$ws.Range('D20').Value = '7.35'
$ws.Range('E20').Value = '  -6.15%  '
$ws.Range('D21').Value = '344.36'
$ws.Range('E21').Value = '  -1.34%  '
$ws.Range('D22').Value = '3.86'
$ws.Range('E22').Value = '  -1.09%  '
$ws.Range('D23').Value = '1.92'
$ws.Range('E23').Value = '  -2.67%  '
$ws.Range('E24').Value = '  -0.13%  '
$ws.Range('D25').Value = '69.57'
$ws.Range('E25').Value = '  -0.84%  '
$ws.Range('D26').Value = '3.89'
$ws.Range('E26').Value = '  -2.67%  '
$ws.Range('D27').Value = '2.618.12'
$ws.Range('E27').Value = '  -2.83%  '
$ws.Range('D28').Value = '8.65'
$ws.Range('E28').Value = '  -1.71%  '
$ws.Range('D29').Value = '0.998'
$ws.Range('E29').Value = '  -0.16%  '
$ws.Range('D30').Value = '0.0₃0872'
$ws.Range('E30').Value = '  -3.14%  '
$ws.Range('D31').Value = '7.67'
$ws.Range('E31').Value = '  -2.40%  '
$ws.Range('D32').Value = '442.21'
$ws.Range('E32').Value = '  -4.29%  '
$ws.Range('B33').Value = 'Fetch.AI'
$ws.Range('C33').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D33').Value = '1.18'
$ws.Range('E33').Value = '  -5.53%  '
$ws.Range('B34').Value = 'FirstDigitalUSD'
$ws.Range('C34').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D34').Value = '0.999'
$ws.Range('E34').Value = '  -0.13%  '
$ws.Range('E35').Value = '  -1.92%  '
$ws.Range('D36').Value = '154.78'
$ws.Range('E36').Value = '  -0.73%  '
$ws.Range('D37').Value = '0.113'
$ws.Range('E37').Value = '  -3.80%  '
$ws.Range('E38').Value = '  +0.00%  '
$ws.Range('D39').Value = '18.12'
$ws.Range('E39').Value = '  -2.79%  '
$ws.Range('E40').Value = '  +0.04%  '
$ws.Range('D41').Value = '0.314'
$ws.Range('E41').Value = '  -0.94%  '
$ws.Range('B42').Value = 'RenderToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D42').Value = '4.58'
$ws.Range('E42').Value = '  -3.70%  '
$ws.Range('B43').Value = 'Stacks'
$ws.Range('C43').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D43').Value = '1.58'
$ws.Range('E43').Value = '  -1.56%  '
$ws.Range('D44').Value = '2.14'
$ws.Range('E44').Value = '  -6.63%  '
$ws.Range('D45').Value = '1.07'
$ws.Range('E45').Value = '  -7.02%  '
$ws.Range('D46').Value = '138.14'
$ws.Range('E46').Value = '  -3.16%  '
$ws.Range('D47').Value = '3.43'
$ws.Range('E47').Value = '  -1.50%  '
$ws.Range('D48').Value = '0.510'
$ws.Range('E48').Value = '  -3.08%  '
$ws.Range('D49').Value = '0.0724'
$ws.Range('E49').Value = '  -0.35%  '
$ws.Range('D50').Value = '0.571'
$ws.Range('E50').Value = '  -1.03%  '
$ws.Range('D51').Value = '0.0920'
$ws.Range('E51').Value = '  -0.96%  '
